$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# =========================================================================
# Add a "Buyer Table" block (rows 10-17), mirroring the "Farmer Table"
# block (rows 1-8), with the id column relabelled "Buyer Id".
# =========================================================================

# Title row (merged, bold, centered, bordered - same look as A1:D1)
$ws.Range("A10").Value = "Buyer Table"
$ws.Range("A10:D10").Merge() | Out-Null
$ws.Range("A10:D10").Borders.Color = 0
$ws.Range("A10:D10").Borders.LineStyle = 1
$ws.Range("A10:D10").Font.Bold = $true
$ws.Range("A10:D10").HorizontalAlignment = -4108

# Header row + 6 data rows (Property Name / Data Type / Key / Default, etc.)
$ws.Range("A2:D8").Copy($ws.Range("A11"))
$ws.Range("A12").Value = "Buyer Id"

# =========================================================================
# Add an "Admin Table" block (rows 19-26), mirroring the "Farmer Table"
# block (rows 1-8), with the id column relabelled "Admin Id".
# =========================================================================

$ws.Range("A19").Value = "Admin Table"
$ws.Range("A19:D19").Merge() | Out-Null
$ws.Range("A19:D19").Borders.Color = 0
$ws.Range("A19:D19").Borders.LineStyle = 1
$ws.Range("A19:D19").Font.Bold = $true
$ws.Range("A19:D19").HorizontalAlignment = -4108

$ws.Range("A2:D8").Copy($ws.Range("A20"))
$ws.Range("A21").Value = "Admin Id"

# =========================================================================
# Update the window selection to match the final state of the workbook.
# =========================================================================
$ws.Range("A29").Select() | Out-Null
